{"js": "// The document contains a single 20-row x 5-column table of simple\n// arithmetic expressions (e.g. \"98-28=\"). The commit replaces the\n// expression text in each of the 100 cells, in row-major order, while\n// leaving every other aspect of the document (formatting, paragraph\n// properties, run properties, table structure, etc.) untouched.\nconst newValues = [\"9+33=\", \"83-20=\", \"51+32=\", \"75-40=\", \"3+1=\", \"48-39=\", \"50-14=\", \"64+32=\", \"16+47=\", \"0+86=\", \"32+60=\", \"2+77=\", \"92-40=\", \"51-48=\", \"79-75=\", \"65-58=\", \"75-15=\", \"42+0=\", \"71+9=\", \"75+8=\", \"40+29=\", \"58-10=\", \"28+22=\", \"55-10=\", \"47+7=\", \"27+20=\", \"14+12=\", \"5+91=\", \"68-12=\", \"23+57=\", \"34+27=\", \"27+43=\", \"53+9=\", \"47-35=\", \"61-48=\", \"56-16=\", \"68+27=\", \"15+42=\", \"42+25=\", \"96-47=\", \"65+12=\", \"57-17=\", \"44-9=\", \"79-50=\", \"33+39=\", \"16+15=\", \"99-84=\", \"26-13=\", \"57+42=\", \"60-26=\", \"33+4=\", \"46+44=\", \"6+71=\", \"53-20=\", \"72-59=\", \"3+51=\", \"52+29=\", \"22+69=\", \"18+61=\", \"61+17=\", \"7+27=\", \"73+0=\", \"87-40=\", \"56+29=\", \"12+43=\", \"32-11=\", \"21+12=\", \"15+1=\", \"56+40=\", \"72-28=\", \"13+4=\", \"28+65=\", \"36+63=\", \"91-53=\", \"78-40=\", \"54-27=\", \"0+58=\", \"41-2=\", \"70-69=\", \"93-76=\", \"87-71=\", \"60+12=\", \"20+13=\", \"44+9=\", \"72-26=\", \"18+20=\", \"19+0=\", \"79+17=\", \"60+23=\", \"19+32=\", \"18+45=\", \"53-30=\", \"54-48=\", \"80-41=\", \"93-56=\", \"96+1=\", \"10+70=\", \"65-50=\", \"85+11=\", \"89-29=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = newValues.length / rowCount;\n\nlet k = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(newValues[k], Word.InsertLocation.replace);\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of simple\n# arithmetic expressions (e.g. \"98-28=\"). The commit replaces the\n# expression text in each of the 100 cells, in row-major order, while\n# leaving every other aspect of the document (formatting, paragraph\n# properties, run properties, table structure, etc.) untouched.\n$newValues = @(\"9+33=\", \"83-20=\", \"51+32=\", \"75-40=\", \"3+1=\", \"48-39=\", \"50-14=\", \"64+32=\", \"16+47=\", \"0+86=\", \"32+60=\", \"2+77=\", \"92-40=\", \"51-48=\", \"79-75=\", \"65-58=\", \"75-15=\", \"42+0=\", \"71+9=\", \"75+8=\", \"40+29=\", \"58-10=\", \"28+22=\", \"55-10=\", \"47+7=\", \"27+20=\", \"14+12=\", \"5+91=\", \"68-12=\", \"23+57=\", \"34+27=\", \"27+43=\", \"53+9=\", \"47-35=\", \"61-48=\", \"56-16=\", \"68+27=\", \"15+42=\", \"42+25=\", \"96-47=\", \"65+12=\", \"57-17=\", \"44-9=\", \"79-50=\", \"33+39=\", \"16+15=\", \"99-84=\", \"26-13=\", \"57+42=\", \"60-26=\", \"33+4=\", \"46+44=\", \"6+71=\", \"53-20=\", \"72-59=\", \"3+51=\", \"52+29=\", \"22+69=\", \"18+61=\", \"61+17=\", \"7+27=\", \"73+0=\", \"87-40=\", \"56+29=\", \"12+43=\", \"32-11=\", \"21+12=\", \"15+1=\", \"56+40=\", \"72-28=\", \"13+4=\", \"28+65=\", \"36+63=\", \"91-53=\", \"78-40=\", \"54-27=\", \"0+58=\", \"41-2=\", \"70-69=\", \"93-76=\", \"87-71=\", \"60+12=\", \"20+13=\", \"44+9=\", \"72-26=\", \"18+20=\", \"19+0=\", \"79+17=\", \"60+23=\", \"19+32=\", \"18+45=\", \"53-30=\", \"54-48=\", \"80-41=\", \"93-56=\", \"96+1=\", \"10+70=\", \"65-50=\", \"85+11=\", \"89-29=\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$k]\n        $k++\n    }\n}\n"}
